$wb = $excel.ActiveWorkbook

# --- CHaMP_Winter_Chinook: fix garbled metric name in B6 ---
$wsChinook = $wb.Worksheets.Item("CHaMP_Winter_Chinook")
$wsChinook.Range("B6").Value = "FishCovSome"
$wsChinook.Range("B6").Select()

# --- CHaMP_Winter_Steelhead: remove duplicate/erroneous "FishCovAqVeg" row (row 6) ---
$wsSteelhead = $wb.Worksheets.Item("CHaMP_Winter_Steelhead")
$wsSteelhead.Activate()
$wsSteelhead.Rows.Item(6).Delete()

# Normalize row-level formatting flags left over from the shift (rows no longer need
# an explicit per-row style once re-numbered), while keeping the centered cell style.
for ($r = 1; $r -le 27; $r++) {
    $wsSteelhead.Rows.Item($r).ClearFormats()
}
$wsSteelhead.Range("A1:E27").HorizontalAlignment = -4108

$wsSteelhead.Range("A6:XFD6").Select()
